$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value2 = 'name'
$ws.Cells.Item(1,2).Value2 = 'address'
$ws.Cells.Item(1,3).Value2 = 'phone'
$ws.Cells.Item(1,4).Value2 = 'description'
$ws.Cells.Item(1,5).Value2 = 'rating'

$ws.Cells.Item(2,1).Value2 = 'Action Pursuit Games of Brandon LLC.'
$ws.Cells.Item(2,2).Value2 = '928 Old Natchez Trace RoadCanton, MS 39046'
$ws.Cells.Item(2,3).Value2 = '(601) 825-1052'
$ws.Cells.Item(2,4).Value2 = 'The premier paintball field in Mississippi. 50 Acres playing fields, HPA only, 150 sets rental equipment. Mini pro shop on site and .50 cal rentals for kids.'
$ws.Cells.Item(2,5).Value2 = '98px'
$ws.Cells.Item(3,1).Value2 = 'Adrenalin Paintball, LLC'
$ws.Cells.Item(3,2).Value2 = '301 Valley DrivePerry, GA 31069'
$ws.Cells.Item(3,3).Value2 = '(478) 987-6055'
$ws.Cells.Item(3,4).Value2 = 'Middle Georgia''s Full Service Paintball proshop for the rec player to the speedball,  tournament player: we have it all!  Planet Eclipse certified Master tech service center for all Planet Eclipse markers. We are a premiere proshop for Big Indian Paintball field.'
$ws.Cells.Item(3,5).Value2 = '50px'
$ws.Cells.Item(4,1).Value2 = 'All Prime Paintball'
$ws.Cells.Item(4,2).Value2 = '1000 bettis academy rdgraniteville, SC 29829'
$ws.Cells.Item(4,3).Value2 = $null
$ws.Cells.Item(4,4).Value2 = $null
$ws.Cells.Item(4,5).Value2 = '100px'
$ws.Cells.Item(5,1).Value2 = 'Arkenstone Paintball'
$ws.Cells.Item(5,2).Value2 = '3198 Cedarcrest RoadAcworth, GA 30101'
$ws.Cells.Item(5,3).Value2 = '(770) 974-2535'
$ws.Cells.Item(5,4).Value2 = 'This is the new Arkenstone Paintball. We are a fast paced woodsball field. Sitting on about ten acres we offer a western town, a castle and many other type fields to play. Come and check us out.'
$ws.Cells.Item(5,5).Value2 = '90px'
$ws.Cells.Item(6,1).Value2 = 'Athens Paintball'
$ws.Cells.Item(6,2).Value2 = '305 Lakeview DriveWinterville, GA 30683'
$ws.Cells.Item(6,3).Value2 = '(706) 583-8600'
$ws.Cells.Item(6,4).Value2 = 'Pro-Shop and 2 Scenario Fields and 2 Town Fields.Walk-On paintball Saturday and Sunday from 10am to 6pm/12pm to 6pm respectively.Private parties 7 days a week. Business, school, church, friends, bachelor and birthdays.Only 10 minutes from Downtown Athens and UGA.'
$ws.Cells.Item(6,5).Value2 = '70px'
$ws.Cells.Item(7,1).Value2 = 'Augusta Paintball'
$ws.Cells.Item(7,2).Value2 = '1818 Barton Chapel RdAugusta, GA 30909'
$ws.Cells.Item(7,3).Value2 = '(256) 267-6244'
$ws.Cells.Item(7,4).Value2 = $null
$ws.Cells.Item(7,5).Value2 = $null
$ws.Cells.Item(8,1).Value2 = 'Battlefield Ridge'
$ws.Cells.Item(8,2).Value2 = '368 County road 8021renzi, MS 38865'
$ws.Cells.Item(8,3).Value2 = '(662) 349-0585'
$ws.Cells.Item(8,4).Value2 = $null
$ws.Cells.Item(8,5).Value2 = $null
$ws.Cells.Item(9,1).Value2 = 'Battlefront Memphis'
$ws.Cells.Item(9,2).Value2 = '6791 Walnut Grove RdMemphis, TN 38115'
$ws.Cells.Item(9,3).Value2 = '(855) 901-7529'
$ws.Cells.Item(9,4).Value2 = 'Outdoor Field on 40 acres of wooded and non-wooded land that offers Outdoor Lasertag, Paintball, Low Impact Paintball, and Arrow Tag! Great prices! Friendly Staff!'
$ws.Cells.Item(9,5).Value2 = $null
$ws.Cells.Item(10,1).Value2 = 'BattleZone Paintball'
$ws.Cells.Item(10,2).Value2 = '3477 King Hall Mill RdBowman, GA 30624'
$ws.Cells.Item(10,3).Value2 = '(706) 614-5834'
$ws.Cells.Item(10,4).Value2 = 'Battle Zone Paintball is a large 165 acre field. The proshop has everything you need and rental packages are available. We are located beside the sandbar and you can rent cabins, kayak, eat at the restaurant and much more. The field ranges from tall hill tops, deep river beds and river flats.'
$ws.Cells.Item(10,5).Value2 = '98px'
$ws.Cells.Item(11,1).Value2 = 'Bear Paw Army Navy Store'
$ws.Cells.Item(11,2).Value2 = '1424 Highway 16 WGriffin, GA 30223'
$ws.Cells.Item(11,3).Value2 = '(770) 228-6634'
$ws.Cells.Item(11,4).Value2 = 'Paintball Markers and Supplies Plus All Your Camouflage Needs'
$ws.Cells.Item(11,5).Value2 = $null
$ws.Cells.Item(12,1).Value2 = 'Big Indian Paintball'
$ws.Cells.Item(12,2).Value2 = '301 Valley DrivePerry, GA 31069'
$ws.Cells.Item(12,3).Value2 = '(478) 951-1331'
$ws.Cells.Item(12,4).Value2 = 'Middle Georgia Largest Paintball Complex with 2 PSP size fields (200ft x 150ft netted) with pits and air in each, a 525ft x 300ft town and 80 ac of woodsball.'
$ws.Cells.Item(12,5).Value2 = '92px'
$ws.Cells.Item(13,1).Value2 = 'Blue Ballz Paintball Supplies'
$ws.Cells.Item(13,2).Value2 = '319 East Thompson StreetThomaston, GA 30256'
$ws.Cells.Item(13,3).Value2 = '(770) 584-4469'
$ws.Cells.Item(13,4).Value2 = 'Your paintball supplier!  Unbeatable prices and sales. Open Friday and Saturday Until 5 pm.'
$ws.Cells.Item(13,5).Value2 = $null
$ws.Cells.Item(14,1).Value2 = 'Boltenwood Paintball'
$ws.Cells.Item(14,2).Value2 = 'Hwy 59 NLinton, IN 47441'
$ws.Cells.Item(14,3).Value2 = '(812) 847-4480'
$ws.Cells.Item(14,4).Value2 = 'Paintball'
$ws.Cells.Item(14,5).Value2 = '80px'
$ws.Cells.Item(15,1).Value2 = 'Brothers In Arms Paintball'
$ws.Cells.Item(15,2).Value2 = '4197 West Highway 84Enterprise, AL 36330'
$ws.Cells.Item(15,3).Value2 = $null
$ws.Cells.Item(15,4).Value2 = 'Speedball field'
$ws.Cells.Item(15,5).Value2 = $null
$ws.Cells.Item(16,1).Value2 = 'Carolina Paintball Fields 7 supplies'
$ws.Cells.Item(16,2).Value2 = '600 Starr Ridge RoadGaffney, SC 29340'
$ws.Cells.Item(16,3).Value2 = '(864) 812-0276'
$ws.Cells.Item(16,4).Value2 = '3 speedball 1 wooded'
$ws.Cells.Item(16,5).Value2 = $null
$ws.Cells.Item(17,1).Value2 = 'Central Alabama Paintball Park'
$ws.Cells.Item(17,2).Value2 = 'Rosebud LaneCalera, AL 35040'
$ws.Cells.Item(17,3).Value2 = '(205) 837-4177'
$ws.Cells.Item(17,4).Value2 = 'A family friendly atmosphere featuring five fields with something different for everybody, including a hyperball field, town field, spools field, and two airball fields. Come enjoy your self at Alabama''s Premier paintball facility!'
$ws.Cells.Item(17,5).Value2 = '92px'
$ws.Cells.Item(18,1).Value2 = 'CHF Paintball'
$ws.Cells.Item(18,2).Value2 = '008 Love RoadHernando, MS 38632'
$ws.Cells.Item(18,3).Value2 = '(662) 429-8787'
$ws.Cells.Item(18,4).Value2 = 'We have woodsball and scenario fields.  An air conditioned/heated pro shop with covered staging areas and a real bathroom.  Cheap prices and good service.'
$ws.Cells.Item(18,5).Value2 = '92px'
$ws.Cells.Item(19,1).Value2 = 'Circle Bar G Paintball Pro Shop'
$ws.Cells.Item(19,2).Value2 = '17 CR 287(Behind Lafayette Spring water tower)Oxford, MS 38655'
$ws.Cells.Item(19,3).Value2 = '(662) 380-3822'
$ws.Cells.Item(19,4).Value2 = 'Paintball Marker and Cylinder Repair Shop. Inspection and repair of cylinders. Teaches Cylinder HAZMAT and Fill Station Compliance IAW with DOT & 49 CFR. Professional Cylinder Cursory Inspection. Certified Advanced Airsmith (PTI) and Professional Cylinder Inspector/Valve Repair Tech (PSI/PCI).'
$ws.Cells.Item(19,5).Value2 = $null
$ws.Cells.Item(20,1).Value2 = 'City Limits Paintball'
$ws.Cells.Item(20,2).Value2 = '1303 Summerlin LaneBastrop, LA 71220'
$ws.Cells.Item(20,3).Value2 = '(318) 283-9494'
$ws.Cells.Item(20,4).Value2 = 'Paintball, paintball, & more paintball!'
$ws.Cells.Item(20,5).Value2 = $null
$ws.Cells.Item(21,1).Value2 = 'Classic Paintball'
$ws.Cells.Item(21,2).Value2 = '1320 Blairs Bridge RoadLithia Springs, GA 30122'
$ws.Cells.Item(21,3).Value2 = '(770) 732-1110'
$ws.Cells.Item(21,4).Value2 = 'Our mission is to provide a Paintball Pro Shop with great products while still being at reasonable prices. We are looking forward to showing you our level Professional Fields.'
$ws.Cells.Item(21,5).Value2 = '92px'
$ws.Cells.Item(22,1).Value2 = 'Co-Op Paintball'
$ws.Cells.Item(22,2).Value2 = '6144 Cumming HighwayBuford, GA 30518'
$ws.Cells.Item(22,3).Value2 = '(404) 510-4937'
$ws.Cells.Item(22,4).Value2 = 'Fun affordable and local to 3 Cities. Co-Op Paintball serves both the new players as well as the Vet!  With rental packages starting at $25 per person w/ everything you need to play, we are comfortable saying you won''t find a better deal anywhere else!'
$ws.Cells.Item(22,5).Value2 = $null
$ws.Cells.Item(23,1).Value2 = 'Code Red Paintball, LLC'
$ws.Cells.Item(23,2).Value2 = '1541 Hwy 35 N.Carthage, MS 39051'
$ws.Cells.Item(23,3).Value2 = '(601) 813-7227'
$ws.Cells.Item(23,4).Value2 = 'When searching for the perfect place to host your paintball event, be it a birthday party, youth/church group, corporate bonding, police/fireman/military training, or just group of friends/family getting together, it is important to have the best equipment and fields.'
$ws.Cells.Item(23,5).Value2 = '100px'
$ws.Cells.Item(24,1).Value2 = 'Cornerstone Paintball'
$ws.Cells.Item(24,2).Value2 = '556 Hwy. 19/41Hampton, GA 30228'
$ws.Cells.Item(24,3).Value2 = '(770) 707-0702'
$ws.Cells.Item(24,4).Value2 = 'Fields and store'
$ws.Cells.Item(24,5).Value2 = $null
$ws.Cells.Item(25,1).Value2 = 'Dan''s Land'
$ws.Cells.Item(25,2).Value2 = 'Harkness RoadJackson, GA 30233'
$ws.Cells.Item(25,3).Value2 = '(770) 851-2000'
$ws.Cells.Item(25,4).Value2 = 'Multiple Fields including X-ball, Hyperball, Speedball, Woodsball, Scenario, and Forts for Attack and Defend'
$ws.Cells.Item(25,5).Value2 = '88px'
$ws.Cells.Item(26,1).Value2 = 'Destruction Paintball'
$ws.Cells.Item(26,2).Value2 = 'Treetop TrailColumbus, MS 39705'
$ws.Cells.Item(26,3).Value2 = '(662) 425-2587'
$ws.Cells.Item(26,4).Value2 = 'Destruction Paintball, in Columbus, Mississippi, is the only full Service Paintball Park in the entire area with Multiple Fields of Play, Rentals, CO2, Compressed Air, a Proshop and Snack Bar.'
$ws.Cells.Item(26,5).Value2 = '100px'
$ws.Cells.Item(27,1).Value2 = 'Doro Sports'
$ws.Cells.Item(27,2).Value2 = '193 Ben Burton CircleBogart, GA 30622'
$ws.Cells.Item(27,3).Value2 = '(706) 549-4900'
$ws.Cells.Item(27,4).Value2 = $null
$ws.Cells.Item(27,5).Value2 = '60px'
$ws.Cells.Item(28,1).Value2 = 'Dosser Works Paintball'
$ws.Cells.Item(28,2).Value2 = '125 Milton AvenueAtlanta, GA 30315'
$ws.Cells.Item(28,3).Value2 = '(678) 927-9359'
$ws.Cells.Item(28,4).Value2 = 'Atlanta''s In-town Paintball. Being the only field inside the perimeter, our view of Atlanta''s skyline gives players a unique urban feel during game play. We have multiple fields, including 2 tournament sized athletic fields. Dosser Works is where the tournament players play.'
$ws.Cells.Item(28,5).Value2 = '90px'
$ws.Cells.Item(29,1).Value2 = 'Eagle Action Sportz'
$ws.Cells.Item(29,2).Value2 = '70 ryner rdHattiesburg, MS 39402'
$ws.Cells.Item(29,3).Value2 = '(601) 341-6735'
$ws.Cells.Item(29,4).Value2 = 'We are a field and pro shop located in Hattiesburg, MS. We offer an amazing experience to both new and veteran players.'
$ws.Cells.Item(29,5).Value2 = $null
$ws.Cells.Item(30,1).Value2 = 'EGS Paintball Field and Supply'
$ws.Cells.Item(30,2).Value2 = '11015 Highway 78Heflin, AL 36264'
$ws.Cells.Item(30,3).Value2 = '(256) 463-4383'
$ws.Cells.Item(30,4).Value2 = '$5.00 Field admission, 20oz CO2 fills $3.00'
$ws.Cells.Item(30,5).Value2 = '20px'
$ws.Cells.Item(31,1).Value2 = 'Embassy Paintball'
$ws.Cells.Item(31,2).Value2 = '100 Crooked Creek RoadAthens, GA 30607'
$ws.Cells.Item(31,3).Value2 = '(706) 215-2387'
$ws.Cells.Item(31,4).Value2 = 'Three exciting fields with the friendliest staff and fun for everyone! CALL TODAY!'
$ws.Cells.Item(31,5).Value2 = $null
$ws.Cells.Item(32,1).Value2 = 'Emerld City Paintball'
$ws.Cells.Item(32,2).Value2 = '1107 D Northlake DriveGreenwood, OR 29649'
$ws.Cells.Item(32,3).Value2 = '(864) 980-9793'
$ws.Cells.Item(32,4).Value2 = $null
$ws.Cells.Item(32,5).Value2 = $null
$ws.Cells.Item(33,1).Value2 = 'Extreme Outdoors'
$ws.Cells.Item(33,2).Value2 = '231 Griggs RoadSix Mile, SC 29682'
$ws.Cells.Item(33,3).Value2 = '(864) 868-4288'
$ws.Cells.Item(33,4).Value2 = 'Woodsball and speedball fields.'
$ws.Cells.Item(33,5).Value2 = $null
$ws.Cells.Item(34,1).Value2 = 'Fast Action Paintball, Inc.'
$ws.Cells.Item(34,2).Value2 = '3331 Crescent RoadBlairsville, GA 30512'
$ws.Cells.Item(34,3).Value2 = $null
$ws.Cells.Item(34,4).Value2 = 'Need some extreme fun? Then visit Fast Action Paintball. We have four outside playing fields, a 1800 sq. ft. showroom with pro shop, snack area, safety area and rest area.'
$ws.Cells.Item(34,5).Value2 = $null
$ws.Cells.Item(35,1).Value2 = 'Firelake Paintball N'' Catfish'
$ws.Cells.Item(35,2).Value2 = '1897 Cheddar RoadAnderson, SC 29627'
$ws.Cells.Item(35,3).Value2 = '(864) 847-4200'
$ws.Cells.Item(35,4).Value2 = 'Firelake has a speedball & woodsball field. There is a proshop with drinks, snacks, great staff, and range of equipment and rental M-98''s. The people are great, the refs are fun. Large staging area.'
$ws.Cells.Item(35,5).Value2 = '100px'
$ws.Cells.Item(36,1).Value2 = 'Fort Benning Paintball Area'
$ws.Cells.Item(36,2).Value2 = 'Dixie RoadFort Benning, GA 31905'
$ws.Cells.Item(36,3).Value2 = '(313) 247-2170'
$ws.Cells.Item(36,4).Value2 = 'A wooded paintball area open to all to have fun and learn more about the sport of Paintball.'
$ws.Cells.Item(36,5).Value2 = '48px'
$ws.Cells.Item(37,1).Value2 = 'GaXtreme Paintball'
$ws.Cells.Item(37,2).Value2 = '105 Cherokee DriveGuyton, GA 31312'
$ws.Cells.Item(37,3).Value2 = '(912) 667-2685'
$ws.Cells.Item(37,4).Value2 = 'Has 2 wooded courses and one speedball field. Very helpful staff and owner.  Always eager to teach new players how to play paintball.'
$ws.Cells.Item(37,5).Value2 = '94px'
$ws.Cells.Item(38,1).Value2 = 'Georgia Paintball'
$ws.Cells.Item(38,2).Value2 = '840 Ernest Barrett ParkwaySuite 144Kennesaw, GA 30144'
$ws.Cells.Item(38,3).Value2 = '(770) 427-2929'
$ws.Cells.Item(38,4).Value2 = 'Atlanta''s Widest Selection of Paintball, Biggest Pro Shop in the Southeast, and more than 50 years of experience, come and visit  us today!'
$ws.Cells.Item(38,5).Value2 = '86px'
$ws.Cells.Item(39,1).Value2 = 'Gorilla War Paintball'
$ws.Cells.Item(39,2).Value2 = '500 Falling Springs RoadRydal, GA 30171'
$ws.Cells.Item(39,3).Value2 = '(770) 757-9152'
$ws.Cells.Item(39,4).Value2 = 'Located off Highway 411, near Cartersville, Fairmount, Adairsville, Canton, Walleska, etc..  Nice field, never too crowded. If you have a tourney coming up, it makes a great field to practice at. It''s cheap and the airball is always open for games.'
$ws.Cells.Item(39,5).Value2 = '68px'
$ws.Cells.Item(40,1).Value2 = 'Highway 441 Paintball'
$ws.Cells.Item(40,2).Value2 = '2856 Highway 441Mountain City, GA 30562'
$ws.Cells.Item(40,3).Value2 = '(706) 982-4994'
$ws.Cells.Item(40,4).Value2 = 'NPS Airball field that is being updated constantly.  Scenario field and more coming in the future.'
$ws.Cells.Item(40,5).Value2 = '100px'
$ws.Cells.Item(41,1).Value2 = 'Hot Shots Paintball - Windy Hill'
$ws.Cells.Item(41,2).Value2 = '100 A Windy Hill RoadMcDonough, GA 30253'
$ws.Cells.Item(41,3).Value2 = '(404) 392-3153'
$ws.Cells.Item(41,4).Value2 = '1 Regulation size X- ball field with room for 2 more, 1 Hyperball field, and a Scenario field.  Staging area, bleachers, Food and Drinks! Family Environment, Team players welcome!'
$ws.Cells.Item(41,5).Value2 = '84px'
$ws.Cells.Item(42,1).Value2 = 'Hullaballoos'
$ws.Cells.Item(42,2).Value2 = '4452 Hwy. 431 SouthEufaula, AL 36027'
$ws.Cells.Item(42,3).Value2 = '(334) 232-4844'
$ws.Cells.Item(42,4).Value2 = 'Paintball fields, bounce houses, paintball shooting gallery and more!'
$ws.Cells.Item(42,5).Value2 = $null
$ws.Cells.Item(43,1).Value2 = 'Insane Paintball'
$ws.Cells.Item(43,2).Value2 = '1200 Wisdom StreetChattanooga, TN 37406'
$ws.Cells.Item(43,3).Value2 = '(423) 624-2121'
$ws.Cells.Item(43,4).Value2 = 'Turf XBall, Town Field, and Woods fields. Indoor and outdoor staging, indoor restrooms, with a 3200 sq ft retail store on site. Minutes from downtown Chattanooga tourist destinations.'
$ws.Cells.Item(43,5).Value2 = '90px'
$ws.Cells.Item(44,1).Value2 = 'Jammin Paintball'
$ws.Cells.Item(44,2).Value2 = '2370 Belltelephone RoadHazlehurst, GA 31539'
$ws.Cells.Item(44,3).Value2 = '(912) 253-9650'
$ws.Cells.Item(44,4).Value2 = 'Speedball style field, 180 ft long by 100 ft wide.'
$ws.Cells.Item(44,5).Value2 = $null
$ws.Cells.Item(45,1).Value2 = 'Knights Crossing Paintball Field'
$ws.Cells.Item(45,2).Value2 = '260 Cowart RoadCommerce, GA 30530'
$ws.Cells.Item(45,3).Value2 = '(770) 530-3984'
$ws.Cells.Item(45,4).Value2 = '12 Acres of tactical paintball which has woods, trench warfare, large field with forts and bunkers with a very nice combination of other items.  Come check us out!  We have been around for over 5 years!'
$ws.Cells.Item(45,5).Value2 = $null
$ws.Cells.Item(46,1).Value2 = 'Knights Crossing Paintball Store'
$ws.Cells.Item(46,2).Value2 = '3701 Atlanta HighwaySuite 20Bogart, GA 30622'
$ws.Cells.Item(46,3).Value2 = '(770) 530-3984'
$ws.Cells.Item(46,4).Value2 = 'Retail Store with paint, air fills, masks, markers, and all paintball accesories, if we don''t have it in stock we can get it next day if needed.'
$ws.Cells.Item(46,5).Value2 = '50px'
$ws.Cells.Item(47,1).Value2 = 'Liberty Paintball LLC'
$ws.Cells.Item(47,2).Value2 = '935 Liberty RoadDanville, AL 35619'
$ws.Cells.Item(47,3).Value2 = '(256) 773-8922'
$ws.Cells.Item(47,4).Value2 = 'Speedball, Woodsball, and proshop! We have certified refs and airsmiths on hand. Food available to all players.'
$ws.Cells.Item(47,5).Value2 = '94px'
$ws.Cells.Item(48,1).Value2 = 'Low Country Paintball'
$ws.Cells.Item(48,2).Value2 = 'Route 1 Box 216Ludowici, GA 31316'
$ws.Cells.Item(48,3).Value2 = '(912) 545-2369'
$ws.Cells.Item(48,4).Value2 = 'X-ball, Sup Air ball, 70+ acres of paintball field, Tournaments & Scenarios. Electro rentals available'
$ws.Cells.Item(48,5).Value2 = '94px'
$ws.Cells.Item(49,1).Value2 = 'Maranatha Camp and Conference Center'
$ws.Cells.Item(49,2).Value2 = '1091 Jeffery RdScottsboro, AL 35769'
$ws.Cells.Item(49,3).Value2 = '(256) 609-3923'
$ws.Cells.Item(49,4).Value2 = 'Camp Maranatha is part of the North Alabama Presbytery. We have a Hyperball and two woods fields. We''re constantly adding to and developing our fields. Please call or email for group rates!'
$ws.Cells.Item(49,5).Value2 = $null
$ws.Cells.Item(50,1).Value2 = 'Maranatha Camp and Conference Center'
$ws.Cells.Item(50,2).Value2 = '1091 Jeffery DriveScottsboro, AL 35769'
$ws.Cells.Item(50,3).Value2 = '(256) 574-4539'
$ws.Cells.Item(50,4).Value2 = 'We are a church camp that has 3 paintball courses right on the beautiful Tennessee River.'
$ws.Cells.Item(50,5).Value2 = $null
$ws.Cells.Item(51,1).Value2 = 'Mississippi Splatter Games'
$ws.Cells.Item(51,2).Value2 = 'Maben Bell Schoolhouse RoadMaben, MS 39750'
$ws.Cells.Item(51,3).Value2 = '(662) 263-4445'
$ws.Cells.Item(51,4).Value2 = '5 to 7 man speedball field composed of Sattelite Dish Bunkers, HyperBall Pipes, and Wooden Bunkers. Scenario field with large barn in the center and woods surrounding the field with a small town of bunkers on each side of the barn.'
$ws.Cells.Item(51,5).Value2 = '86px'
$ws.Cells.Item(52,1).Value2 = 'Mississippi State Paintball'
$ws.Cells.Item(52,2).Value2 = 'Recplex Intramural Baseball & SoftballStarkville, MS 39759'
$ws.Cells.Item(52,3).Value2 = '(662) 769-2263'
$ws.Cells.Item(52,4).Value2 = 'Speedball field, great for 3v3 but can accommodate 5v5. It is composed of inflatable bunkers and permanent bunkers.'
$ws.Cells.Item(52,5).Value2 = '70px'
$ws.Cells.Item(53,1).Value2 = 'Mobile Wargamez'
$ws.Cells.Item(53,2).Value2 = '504 choctaw rdenterprise, AL 36330'
$ws.Cells.Item(53,3).Value2 = '(334) 494-9674'
$ws.Cells.Item(53,4).Value2 = 'War Gamez, LLC provides Enterprise and surrounding areas with the first "Mobile Paintball" rentals.  We will bring the guns, air, and paint to you.  All you have to do is to pick a battlefield. If you don''t have a battlefield, we can also provide that!  But just think, of the advantage you would have on your own "Turf."'
$ws.Cells.Item(53,5).Value2 = $null
$ws.Cells.Item(54,1).Value2 = 'Monroe Paintball Works'
$ws.Cells.Item(54,2).Value2 = '320 South Madison AvenueMonroe, GA 30655'
$ws.Cells.Item(54,3).Value2 = '(770) 881-7421'
$ws.Cells.Item(54,4).Value2 = '2 Indoor Fields1 Outdoor Field(at the present time)Astroturf InsideFully Stocked Pro-shopExperienced Service Dept.Cheap Paint, Supplies, Field Fees, etc.770-881-7421Call for Info on BYOP.'
$ws.Cells.Item(54,5).Value2 = $null
$ws.Cells.Item(55,1).Value2 = 'Montgomery Paintball'
$ws.Cells.Item(55,2).Value2 = '5765 Carmichael PkwyMontgomery, AL 36117'
$ws.Cells.Item(55,3).Value2 = '(334) 657-7878'
$ws.Cells.Item(55,4).Value2 = $null
$ws.Cells.Item(55,5).Value2 = $null
$ws.Cells.Item(56,1).Value2 = 'Moxie Paintball LLC'
$ws.Cells.Item(56,2).Value2 = '7519 Fortson RoadSuite 1-CColumbus, GA 31909'
$ws.Cells.Item(56,3).Value2 = '(706) 221-6700'
$ws.Cells.Item(56,4).Value2 = 'Moxie Paintball is Columbus'' only full-service paintball store. Specializing in a full-line of paintball products and brands, available for purchase or special order. With a factory-authorized technician on site, Moxie can repair or modify your paintball marker.'
$ws.Cells.Item(56,5).Value2 = '100px'
$ws.Cells.Item(57,1).Value2 = 'Mt. Doom Paintball Field'
$ws.Cells.Item(57,2).Value2 = '3071 County Road 515Hanceville, AL 35077'
$ws.Cells.Item(57,3).Value2 = '(256) 339-1601'
$ws.Cells.Item(57,4).Value2 = 'Oldest field in the southeastern USA.  Full rentals for 60 players with Woodsball, Speedball, and Scenario Fields.  200 acres for your enjoyment.  Since 1987 we have been serving the paintball sport.'
$ws.Cells.Item(57,5).Value2 = '98px'
$ws.Cells.Item(58,1).Value2 = 'New Gen Paintball'
$ws.Cells.Item(58,2).Value2 = '108 Semper Fidelis RoadEasley, SC 29640'
$ws.Cells.Item(58,3).Value2 = '(864) 859-5561'
$ws.Cells.Item(58,4).Value2 = 'South Carolina''s largest paintball facility.  40,000 square foot building on site with Pro Shop, Party Rooms, Showers, Team Rooms, Skate Park and more.  5 fields on site.  2 airball fields, 1 hyperball field, 1 concept field and 1 wooded scenario field.'
$ws.Cells.Item(58,5).Value2 = '88px'
$ws.Cells.Item(59,1).Value2 = 'Nitro Paintball'
$ws.Cells.Item(59,2).Value2 = '204a Lower Bethany RoadCanton, GA 30114'
$ws.Cells.Item(59,3).Value2 = '(678) 793-8648'
$ws.Cells.Item(59,4).Value2 = 'We were created to provide a quality play with a friendly environment at prices anyone can afford. At our fields you can have a unique woods ball experience, while preserving the safety of yourself and others. We have a tournament-style speedball field.'
$ws.Cells.Item(59,5).Value2 = '96px'
$ws.Cells.Item(60,1).Value2 = 'North Alabama Paintball'
$ws.Cells.Item(60,2).Value2 = 'Baker LNTuscumbia, AL 35674'
$ws.Cells.Item(60,3).Value2 = $null
$ws.Cells.Item(60,4).Value2 = $null
$ws.Cells.Item(60,5).Value2 = '96px'
$ws.Cells.Item(61,1).Value2 = 'North Side Fun Park'
$ws.Cells.Item(61,2).Value2 = '100 Southern DriveLaGrange, GA 30240'
$ws.Cells.Item(61,3).Value2 = '(706) 443-5221'
$ws.Cells.Item(61,4).Value2 = 'Northside Fun Park is one of Georgia''s only facility with indoor paintball. We also have two outdoor paintball fields featuring different courses and indoor laser tag. Our paintball facility is just 60 minutes from Atlanta, directly south on I-85.'
$ws.Cells.Item(61,5).Value2 = '74px'
$ws.Cells.Item(62,1).Value2 = 'On Target Paintball Ga.'
$ws.Cells.Item(62,2).Value2 = '7450 Hawkinsville RoadMacon, GA 31216'
$ws.Cells.Item(62,3).Value2 = '(478) 714-2003'
$ws.Cells.Item(62,4).Value2 = 'Loads of fun with airball, pipe fields, scenario field and a pro shop with staff on hand at all times.BYOP on Sundays or buy a case and free field fee.'
$ws.Cells.Item(62,5).Value2 = '92px'
$ws.Cells.Item(63,1).Value2 = 'Paintball Atlanta'
$ws.Cells.Item(63,2).Value2 = '5315 Shiloh RoadAlpharetta, GA 30005'
$ws.Cells.Item(63,3).Value2 = '(770) 594-0912'
$ws.Cells.Item(63,4).Value2 = 'The field is one of the oldest and most recognized fields in paintball. Paintball Atlanta has been in continuous operation for over 30 years.New equipment for renters. Diverse field maps for players enjoyment.Paintball Atlanta provides a fun, safe, exciting paintball experience.'
$ws.Cells.Item(63,5).Value2 = '72px'
$ws.Cells.Item(64,1).Value2 = 'Paintball Bunker (The Bunker)'
$ws.Cells.Item(64,2).Value2 = '2424 Mount Pleasant RoadHernando, MS 38632'
$ws.Cells.Item(64,3).Value2 = '(662) 470-4843'
$ws.Cells.Item(64,4).Value2 = 'Paintball and Airsoft Store'
$ws.Cells.Item(64,5).Value2 = '50px'
$ws.Cells.Item(65,1).Value2 = 'Paintball Central'
$ws.Cells.Item(65,2).Value2 = '870 Raymond Hudnell RoadCollinsville, MS 39325'
$ws.Cells.Item(65,3).Value2 = '(601) 917-0849'
$ws.Cells.Item(65,4).Value2 = 'New field & store near Collinsville & Philadelphia off Hwy 19.  Currently has speedball & woodsball courses & more to come.  New equipment arriving weekly in our store as well.  For more info call 601-917-0849 or e-mail rhudnell69@aol.com'
$ws.Cells.Item(65,5).Value2 = '20px'
$ws.Cells.Item(66,1).Value2 = 'Paintball Central'
$ws.Cells.Item(66,2).Value2 = '700 Garlington Rd.Suite IGreenville, SC 29607'
$ws.Cells.Item(66,3).Value2 = '(864) 254-0111'
$ws.Cells.Item(66,4).Value2 = 'Paintball Central, the southeast’s largest paintball retailer, we’ve been providing quality products and excellent service to the Carolina’s for more than 15 years. 3 retail locations and 3 paintball parks in North and South Carolina.'
$ws.Cells.Item(66,5).Value2 = '98px'
$ws.Cells.Item(67,1).Value2 = 'Paintball GI'
$ws.Cells.Item(67,2).Value2 = '4056 Pacolet HighwayPacolet, SC 29372'
$ws.Cells.Item(67,3).Value2 = '(864) 474-4233'
$ws.Cells.Item(67,4).Value2 = 'Great Prices, Fast Shipping, Friendly Customer Service.  For all your paintball needs.'
$ws.Cells.Item(67,5).Value2 = $null
$ws.Cells.Item(68,1).Value2 = 'Paintball is Good'
$ws.Cells.Item(68,2).Value2 = '8651 Serene DrMcCalla, AL 35111'
$ws.Cells.Item(68,3).Value2 = '(205) 477-6067'
$ws.Cells.Item(68,4).Value2 = 'Paintball Field and store'
$ws.Cells.Item(68,5).Value2 = '20px'
$ws.Cells.Item(69,1).Value2 = 'Paintball Statesboro'
$ws.Cells.Item(69,2).Value2 = '912-764-9799Statesboro, GA 30458'
$ws.Cells.Item(69,3).Value2 = '(912) 764-9799'
$ws.Cells.Item(69,4).Value2 = 'Indoor field with 3 different types of courses! Nice astro turf. Brand new equipment.'
$ws.Cells.Item(69,5).Value2 = '70px'
$ws.Cells.Item(70,1).Value2 = 'Paintball-Outfitter'
$ws.Cells.Item(70,2).Value2 = '1180 Lagrange HighwayGreenville, GA 30222'
$ws.Cells.Item(70,3).Value2 = '(706) 882-8721'
$ws.Cells.Item(70,4).Value2 = 'Three fields: an old western town, a speedball field made up of wooden spools and barrels, and a castle.2 scenario fields1 speedball field100 acre woods'
$ws.Cells.Item(70,5).Value2 = '22px'
$ws.Cells.Item(71,1).Value2 = 'Parkers Paintball Field'
$ws.Cells.Item(71,2).Value2 = '2866 Newton Conehatta RoadLawrence, MS 39336'
$ws.Cells.Item(71,3).Value2 = '(601) 683-0966'
$ws.Cells.Item(71,4).Value2 = 'Junk car speed ball field, woods, hyper-ball on the way and air bags field. 40 acres with some of the field under consrtuction as this is a new field.'
$ws.Cells.Item(71,5).Value2 = $null
$ws.Cells.Item(72,1).Value2 = 'PBC Sports Park - Greenville'
$ws.Cells.Item(72,2).Value2 = '300 Tucapau RoadWellford, SC 29385(866) 421-PLAY'
$ws.Cells.Item(72,3).Value2 = $null
$ws.Cells.Item(72,4).Value2 = 'PBC features a turf covered Xball field, a town field, 1 Hyperball Field, 1 Mounds Field, and multiple wooded scenario fields with pillboxes, bunkers, villages, castles, and large forts! Up to 4500 psi air fills and Co2. Walk-on and Private play available.'
$ws.Cells.Item(72,5).Value2 = '90px'
$ws.Cells.Item(73,1).Value2 = 'Phoenix Paintball Division'
$ws.Cells.Item(73,2).Value2 = '24970 Highway 72Athens, AL 35611'
$ws.Cells.Item(73,3).Value2 = '(256) 232-9599'
$ws.Cells.Item(73,4).Value2 = 'Over 15 acres of rec-ball fields.  We rent Tippmann Model 98''s with HPA tanks.  Looking for fast action -- try our lighted tournament airball fields.  Our staff puts their many years of paintball experience to work for you ensuring a safe and fun time.'
$ws.Cells.Item(73,5).Value2 = '66px'
$ws.Cells.Item(74,1).Value2 = 'Planet Paintball'
$ws.Cells.Item(74,2).Value2 = '380 Blackwood Store RoadMoore, SC 29369'
$ws.Cells.Item(74,3).Value2 = '(864) 525-3200'
$ws.Cells.Item(74,4).Value2 = 'Planet Paintball offers 6 different fields. We have airball, hyperball, city scenario field and our woods course. We have HPA and CO2 fills. We rent Tippmans and offer IONS as upgrades. We have a stocked proshop, gun tech on premises and courteous refs.'
$ws.Cells.Item(74,5).Value2 = '94px'
$ws.Cells.Item(75,1).Value2 = 'Raging Tiger Paintball'
$ws.Cells.Item(75,2).Value2 = '261 Grace Chapel RoadEnoree, SC 29335'
$ws.Cells.Item(75,3).Value2 = '(864) 940-1228'
$ws.Cells.Item(75,4).Value2 = 'Raging Tiger Paintball - Scenario Games & Online Store.'
$ws.Cells.Item(75,5).Value2 = $null
$ws.Cells.Item(76,1).Value2 = 'Reality paintball'
$ws.Cells.Item(76,2).Value2 = '2462 Lawrence Cove RoadEva, AL 35621'
$ws.Cells.Item(76,3).Value2 = '(256) 318-5593'
$ws.Cells.Item(76,4).Value2 = 'A small field that has a ghetto speedball field and several woodsball fields- Paint and nitro fills are available- FREE TO PLAY!!!'
$ws.Cells.Item(76,5).Value2 = $null
$ws.Cells.Item(77,1).Value2 = 'Rebel Paintball LLC'
$ws.Cells.Item(77,2).Value2 = '28054 Highway 6Sardis, MS 38666'
$ws.Cells.Item(77,3).Value2 = '(662) 380-0997'
$ws.Cells.Item(77,4).Value2 = '4 fields: Speedball, Mad Max Road Course, Pallet Course, Storm the Hill (woodsball)'
$ws.Cells.Item(77,5).Value2 = '66px'
$ws.Cells.Item(78,1).Value2 = 'Red Fox'
$ws.Cells.Item(78,2).Value2 = '1 Red Fox RunWoodruff, SC 29388'
$ws.Cells.Item(78,3).Value2 = '(864) 386-7304'
$ws.Cells.Item(78,4).Value2 = 'Come and play at Red Fox Games today and see what the South''s LARGEST paintball and airsoft field has to offer you! 162 acres with 13 fields.'
$ws.Cells.Item(78,5).Value2 = '74px'
$ws.Cells.Item(79,1).Value2 = 'Righteous Paintball at Northside Fun Park'
$ws.Cells.Item(79,2).Value2 = '100 Southern DriveLaGrange, GA 30240'
$ws.Cells.Item(79,3).Value2 = '(706) 302-9129'
$ws.Cells.Item(79,4).Value2 = 'An up and coming paintball facility in LaGrange, Ga with things for all ages to do.  We also have Laser Tag.'
$ws.Cells.Item(79,5).Value2 = '100px'
$ws.Cells.Item(80,1).Value2 = 'River Bend Resort'
$ws.Cells.Item(80,2).Value2 = '1000 Wilkie Bridge RoadInman, SC 29349'
$ws.Cells.Item(80,3).Value2 = '(864) 621-8965'
$ws.Cells.Item(80,4).Value2 = 'At River Bend, instruction in paintball, shotgun, pistol, rifle and archery at our Summer Camp. The River Bend Lodge, the heart of our resort, has a pro shop, dining rooms, locker rooms, steam rooms, bar and lounge, fireplace and conference areas.'
$ws.Cells.Item(80,5).Value2 = $null
$ws.Cells.Item(81,1).Value2 = 'Rose Hill War Zone'
$ws.Cells.Item(81,2).Value2 = '33071 Rose Hill RoadDozier, AL 36028'
$ws.Cells.Item(81,3).Value2 = '(334) 465-9733'
$ws.Cells.Item(81,4).Value2 = 'Paintball Field and Store'
$ws.Cells.Item(81,5).Value2 = '26px'
$ws.Cells.Item(82,1).Value2 = 'RWB Paintball'
$ws.Cells.Item(82,2).Value2 = '3001 New Cut RoadInman, SC 29349'
$ws.Cells.Item(82,3).Value2 = '(864) 439-8580'
$ws.Cells.Item(82,4).Value2 = 'Offering all the top name brands of guns and gear.  If we aint got it in stock we can get it!  OPEN SUNDAYS!'
$ws.Cells.Item(82,5).Value2 = $null
$ws.Cells.Item(83,1).Value2 = 'Sand Mountain Shooters Club'
$ws.Cells.Item(83,2).Value2 = '626 Bloodworth RoadBoaz, AL 35956'
$ws.Cells.Item(83,3).Value2 = '(256) 593-8027'
$ws.Cells.Item(83,4).Value2 = '3 shooting ranges, small field, two snakes and other great bunkers.'
$ws.Cells.Item(83,5).Value2 = '100px'
$ws.Cells.Item(84,1).Value2 = 'Shop4paintball'
$ws.Cells.Item(84,2).Value2 = '2801 Wade Hampton Boulevard#224Taylors, SC 29687'
$ws.Cells.Item(84,3).Value2 = '(888) 719-6212'
$ws.Cells.Item(84,4).Value2 = 'Everything paintball, markers, tanks, paint, clothes, barrells, mods'
$ws.Cells.Item(84,5).Value2 = $null
$ws.Cells.Item(85,1).Value2 = 'Sledgehammer Paintball F&PS'
$ws.Cells.Item(85,2).Value2 = '3918 GA. Highway 355Buena Vista, GA 31803'
$ws.Cells.Item(85,3).Value2 = '(229) 649-2122'
$ws.Cells.Item(85,4).Value2 = 'We have 4 great fields to play on, including our New PSP Regulation Speedball Field,  our Western Town: "Hangover", fast playing scenario field: "Ft. Apache" our wooded field with forts and Desert Storm with hills and trenches. Check out our site for pics, prices and other information!!'
$ws.Cells.Item(85,5).Value2 = '44px'
$ws.Cells.Item(86,1).Value2 = 'Southeast Paintball'
$ws.Cells.Item(86,2).Value2 = '10455 Highway 431Roanoke, AL 36274'
$ws.Cells.Item(86,3).Value2 = '(256) 563-5768'
$ws.Cells.Item(86,4).Value2 = 'Southeast paintball is located in Roanoke, AL and currently features five nice wooded fields with an NPPL Sup Air field that should be up sometime this summer. They have reasonable rates and the owners are great people. Definitly check this place out!'
$ws.Cells.Item(86,5).Value2 = $null
$ws.Cells.Item(87,1).Value2 = 'Spirit Cycles'
$ws.Cells.Item(87,2).Value2 = '916 Washington StJefferson, GA 30549'
$ws.Cells.Item(87,3).Value2 = '(706) 367-5050'
$ws.Cells.Item(87,4).Value2 = 'Don''t let the name fool you, we have a full service paintball shop inside!'
$ws.Cells.Item(87,5).Value2 = $null
$ws.Cells.Item(88,1).Value2 = 'Splat Zone Paintball'
$ws.Cells.Item(88,2).Value2 = '1550 Avalon DriveCottondale, AL 35453'
$ws.Cells.Item(88,3).Value2 = $null
$ws.Cells.Item(88,4).Value2 = 'Airball, Lego Speedball, Mini Lego Speedball, Desert Assualt, Woods, King of the Hill.'
$ws.Cells.Item(88,5).Value2 = '90px'
$ws.Cells.Item(89,1).Value2 = 'Splattered Woods Paintball'
$ws.Cells.Item(89,2).Value2 = '2151 Cantelou RoadMontgomery, AL 36108'
$ws.Cells.Item(89,3).Value2 = '(334) 320-7011'
$ws.Cells.Item(89,4).Value2 = 'Home Field of E and S Hobbies, Inc. Woodsball, Speedball, Birthday Parties, Youth Group Events, Team building, Tournaments and Scenario Games.'
$ws.Cells.Item(89,5).Value2 = '96px'
$ws.Cells.Item(90,1).Value2 = 'Splatters Outdoor Adventure'
$ws.Cells.Item(90,2).Value2 = '862 Sparta Highway NEMilledgeville, GA 31061'
$ws.Cells.Item(90,3).Value2 = '(478) 451-0705'
$ws.Cells.Item(90,4).Value2 = 'Splatters Outdoor Adventure is one of the nicest paintball fields around.  It is located in Milledgeville, Georgia and has one air ball field, one speedball field and a woodsball field that is over 8 acres.'
$ws.Cells.Item(90,5).Value2 = $null
$ws.Cells.Item(91,1).Value2 = 'Squirrels Nest Paintball'
$ws.Cells.Item(91,2).Value2 = '129 Collins RdEasley, SC 29642'
$ws.Cells.Item(91,3).Value2 = '(864) 760-8448'
$ws.Cells.Item(91,4).Value2 = 'Tactical Paintball Field / Store'
$ws.Cells.Item(91,5).Value2 = '96px'
$ws.Cells.Item(92,1).Value2 = 'Statesboro Paintball'
$ws.Cells.Item(92,2).Value2 = '158 East Parrish StStatesboro, GA 30458'
$ws.Cells.Item(92,3).Value2 = '(912) 764-9799'
$ws.Cells.Item(92,4).Value2 = 'Indoor Paintball'
$ws.Cells.Item(92,5).Value2 = $null
$ws.Cells.Item(93,1).Value2 = 'Tattoo Paintball'
$ws.Cells.Item(93,2).Value2 = '626 Bloodworth RoadGranada Hills, GA 91344'
$ws.Cells.Item(93,3).Value2 = '(818) 366-3079'
$ws.Cells.Item(93,4).Value2 = $null
$ws.Cells.Item(93,5).Value2 = $null
$ws.Cells.Item(94,1).Value2 = 'Taylor Made Gaming'
$ws.Cells.Item(94,2).Value2 = '11075 Wares Ferry RdMontgomery, AL 36117'
$ws.Cells.Item(94,3).Value2 = '(334) 538-4944'
$ws.Cells.Item(94,4).Value2 = 'We are a new Paintball field located on the property of Jenkins Creek Adventures on Wares Ferry Rd. We welcome all levels of play and anyone in who just wants to come out and have a good time.'
$ws.Cells.Item(94,5).Value2 = $null
$ws.Cells.Item(95,1).Value2 = 'The Hq Army-Navy Stores, Inc.'
$ws.Cells.Item(95,2).Value2 = '1735 Montebello Town CenterMontebello, AL 90640'
$ws.Cells.Item(95,3).Value2 = '(323) 727-9852'
$ws.Cells.Item(95,4).Value2 = $null
$ws.Cells.Item(95,5).Value2 = $null
$ws.Cells.Item(96,1).Value2 = 'The Paintball Store (National Pro PBS)'
$ws.Cells.Item(96,2).Value2 = '46 Keri LaneAthens, GA 30607'
$ws.Cells.Item(96,3).Value2 = '(706) 552-1225'
$ws.Cells.Item(96,4).Value2 = 'The Paintball Store is locally owned with a wide variety of paintball guns, paintballs, and other accessories at the best prices guaranteed, we have woodball and speedball fields and we can help you find the paintball equipment that you need.'
$ws.Cells.Item(96,5).Value2 = $null
$ws.Cells.Item(97,1).Value2 = 'War With Friends Paintball'
$ws.Cells.Item(97,2).Value2 = '745 Morris RdHiram, GA 30127'
$ws.Cells.Item(97,3).Value2 = '(404) 771-6626'
$ws.Cells.Item(97,4).Value2 = 'Battle War Gaming Paintball club. Free to play, combines padded melee weapons, throwing weapons, shields, armored vests, with single shot paintball guns and slingshots. CTF games hosted twice a month. Warriors and Barbarian only need apply'
$ws.Cells.Item(97,5).Value2 = $null
$ws.Cells.Item(98,1).Value2 = 'Weekend Warrior Sportz'
$ws.Cells.Item(98,2).Value2 = '107 Old Laurens RoadSimpsonville, SC 29681'
$ws.Cells.Item(98,3).Value2 = '(864) 688-0123'
$ws.Cells.Item(98,4).Value2 = 'Custom Safety Equipment such as barrel covers, dead rags, safety signs & netting. Our customer service and product knowledge is unparalleled, we believe in happy customers, not just making sales for things you don''t need.'
$ws.Cells.Item(98,5).Value2 = $null
$ws.Cells.Item(99,1).Value2 = 'Whitetail Ridge'
$ws.Cells.Item(99,2).Value2 = '444 Birmingham Ridge RoadBlue Springs, MS 38828'
$ws.Cells.Item(99,3).Value2 = '(662) 869-2925'
$ws.Cells.Item(99,4).Value2 = 'Large Club House, Full Service Pro Shop, Air-Ball, Hyper-Ball & Wooded Fields, Skeet, Trap, 5 Stand, Sporting Clays, Rifle & Pistol Ranges, Archery, & Fishing'
$ws.Cells.Item(99,5).Value2 = '96px'
$ws.Cells.Item(100,1).Value2 = 'Wildfire Paintball  - Conyers'
$ws.Cells.Item(100,2).Value2 = '2051 Highway 138 NEConyers, GA 30013'
$ws.Cells.Item(100,3).Value2 = '(770) 817-0521'
$ws.Cells.Item(100,4).Value2 = 'Wildfire Paintball Games outdoor field located in Conyers, GA.  The field has an  a speedball field, and a woodsball field.'
$ws.Cells.Item(100,5).Value2 = '90px'
$ws.Cells.Item(101,1).Value2 = 'Wildfire Paintball - Snellville (Indoor)'
$ws.Cells.Item(101,2).Value2 = '3725-C Stone Mountain HighwaySnellville, GA 30039'
$ws.Cells.Item(101,3).Value2 = '(770) 982-8180'
$ws.Cells.Item(101,4).Value2 = 'The only indoor field in the state.  We now have a fully astro-turfed X-Ball field with stadium lighting.  Wildfire also has a fully stocked store with an experienced and talented airsmith.  Great place to play when the other fields are closed for rain.'
$ws.Cells.Item(101,5).Value2 = '96px'
$ws.Cells.Item(102,1).Value2 = 'Woodlands Paintball'
$ws.Cells.Item(102,2).Value2 = '1229 Mt. Zion Rd.Hayden, AL 35079'
$ws.Cells.Item(102,3).Value2 = '(205) 712-5363'
$ws.Cells.Item(102,4).Value2 = 'Woodlands Paintball is located in Hayden, Alabama.SCENARIO WOODSBALLTHE BOX: 15 ACRES OF PLAYING FIELDS ALL SIZES OF COURSE AVAILABLE FOR PLAYDRINKS AND SNACKS AVAILABLEWALK-ONS, TEAMS,  AND GROUPS WELCOME'
$ws.Cells.Item(102,5).Value2 = '100px'
$ws.Cells.Item(103,1).Value2 = '907 Paintball'
$ws.Cells.Item(103,2).Value2 = '700 W. Klatt RdLocated inside the South Anchorage Sports ParkAnchorage, AK 99515'
$ws.Cells.Item(103,3).Value2 = '(907) 727-7614'
$ws.Cells.Item(103,4).Value2 = '907 Paintball is the home of the midnight sun tournament series, open to all Alaskan paintballers.'
$ws.Cells.Item(103,5).Value2 = '92px'
$ws.Cells.Item(104,1).Value2 = 'Allegiance Paintball'
$ws.Cells.Item(104,2).Value2 = 'Saprae CreekFort McMurray, AB Canada'
$ws.Cells.Item(104,3).Value2 = '(780) 799-3669'
$ws.Cells.Item(104,4).Value2 = $null
$ws.Cells.Item(104,5).Value2 = $null
$ws.Cells.Item(105,1).Value2 = 'Arctic Fire & Safety'
$ws.Cells.Item(105,2).Value2 = '702 30th AvenueFairbanks, AK 99701'
$ws.Cells.Item(105,3).Value2 = '(907) 452-7806'
$ws.Cells.Item(105,4).Value2 = 'Air and CO2 fills but no longer carries gear.'
$ws.Cells.Item(105,5).Value2 = $null
$ws.Cells.Item(106,1).Value2 = 'Solid Paintball'
$ws.Cells.Item(106,2).Value2 = 'Solid Paintball WayAnchorage, AK 99503'
$ws.Cells.Item(106,3).Value2 = '(907) 830-8770'
$ws.Cells.Item(106,4).Value2 = 'Home of Division II Nppl Super 7 paintball Team Solid.'
$ws.Cells.Item(106,5).Value2 = '20px'
$ws.Cells.Item(107,1).Value2 = 'Xcalibre Paintball Park'
$ws.Cells.Item(107,2).Value2 = '63116  Hwy 831Long Lake, AB Canada'
$ws.Cells.Item(107,3).Value2 = '(780) 576-2211'
$ws.Cells.Item(107,4).Value2 = '2 Incredible Speedball Feilds -Nppl and NXL Regulation size8 other woodsball/rec feilds to choose from'
$ws.Cells.Item(107,5).Value2 = '84px'
$ws.Cells.Item(108,1).Value2 = 'Yellowknife Paintball'
$ws.Cells.Item(108,2).Value2 = '10000 kam lake roadyellowknife, AB x1a3y1Canada'
$ws.Cells.Item(108,3).Value2 = $null
$ws.Cells.Item(108,4).Value2 = 'Single speedball field operating 1-2 times per week'
$ws.Cells.Item(108,5).Value2 = $null

# Column widths (closest achievable given engine pixel snapping)
$ws.Columns.Item(2).ColumnWidth = 72.33333333333333
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666

# Sheet view: selection + zoom
$ws.Activate()
$ws.Range("C1:C1048576").Select()
$excel.ActiveWindow.Zoom = 100
